$d = $word.ActiveDocument

# 1. Insert "There are two cases. " immediately before "Choose the case that is correct."
$d.Content.Find.Execute("Choose the case that is correct.", $false, $false, $false, $false, $false,
                         $true, 1, $false, "There are two cases. Choose the case that is correct.", 2)

# 2. Relocate the "_GoBack" bookmark so that it now sits right after the newly
#    inserted "There are two cases. " text and right before "Choose the case
#    that is correct." (Adding a bookmark with an already-existing name moves it).
$r = $d.Content
$r.Find.Execute("Choose the case that is correct.", $false, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)
$r.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r)

# 3. Rename the UserChoices dictionary keys (two occurrences each, Case 1 & Case 2 sections).
$d.Content.Find.Execute("referenceMeasuredFileName", $true, $false, $false, $false, $false,
                         $true, 1, $false, "referenceFileToTune", 2)
$d.Content.Find.Execute("referenceLiteratureFileName", $true, $false, $false, $false, $false,
                         $true, 1, $false, "referenceFileToMatch", 2)

# 4. Remove the now-obsolete "naming is confusing" remarks (two occurrences,
#    replaced in a single Replace-All style pass).
$d.Content.Find.Execute("The naming of this variable is confusing because it was originally designed for Case 1.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
